# Updated cryptos list on Tue Nov  5 05:32:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value. Only D (Price) and E (Volume(1h)) change.
# Values are written as text to preserve the original "inlineStr" cell type
# (these are display strings, not numeric values, e.g. "68.391.83").

$updates = @{
    2  = @{ D = "68.391.83"; E = "  -0.90%  " }
    3  = @{ D = "2.425.91";  E = "  -1.66%  " }
    4  = @{ E = "  -0.02%  " }
    5  = @{ D = "557.99";    E = "  -0.36%  " }
    6  = @{ D = "160.09";    E = "  -1.27%  " }
    8  = @{ E = "  +0.50%  " }
    9  = @{ E = "  +8.20%  " }
    10 = @{ E = "  -1.64%  " }
    11 = @{ E = "  -0.20%  " }
    12 = @{ E = "  -5.28%  " }
    13 = @{ D = "68.281.36"; E = "  -0.90%  " }
    14 = @{ D = "2.871.37" }
    15 = @{ E = "  +2.97%  " }
    16 = @{ E = "  -2.32%  " }
    17 = @{ D = "2.426.38";  E = "  -1.76%  " }
    18 = @{ D = "10.46";     E = "  -2.42%  " }
    19 = @{ D = "334.59";    E = "  -0.65%  " }
    20 = @{ E = "  -1.45%  " }
    21 = @{ E = "  +0.81%  " }
    22 = @{ D = "1.92";      E = "  +1.89%  " }
    23 = @{ E = "  -0.02%  " }
    24 = @{ D = "66.68";     E = "  -0.44%  " }
    25 = @{ D = "3.67";      E = "  +0.19%  " }
    26 = @{ D = "2.554.00";  E = "  -1.80%  " }
    27 = @{ D = "1.00";      E = "  +0.22%  " }
    28 = @{ D = "8.20";      E = "  -0.13%  " }
    29 = @{ D = ("0.0{0}0816" -f [char]0x2083); E = "  -0.16%  " }
    30 = @{ D = "7.14";      E = "  -0.56%  " }
    31 = @{ D = "0.999";     E = "  +0.01%  " }
    32 = @{ D = "424.60";    E = "  -1.19%  " }
    33 = @{ E = "  +0.29%  " }
    34 = @{ E = "  -0.68%  " }
    35 = @{ D = "159.80";    E = "  +1.89%  " }
    36 = @{ D = "19.03" }
    37 = @{ E = "  -0.02%  " }
    38 = @{ E = "  +0.59%  " }
    39 = @{ E = "  -3.53%  " }
    40 = @{ E = "  -0.93%  " }
    41 = @{ D = "4.33";      E = "  -1.86%  " }
    42 = @{ E = "  +1.52%  " }
    43 = @{ E = "  +0.34%  " }
    44 = @{ D = "2.03";      E = "  -1.14%  " }
    45 = @{ D = "131.58" }
    46 = @{ E = "  -0.34%  " }
    47 = @{ D = "0.0712";    E = "  -0.30%  " }
    48 = @{ E = "  -0.62%  " }
    49 = @{ D = "0.556";     E = "  -0.78%  " }
    50 = @{ D = "0.0914";    E = "  +0.08%  " }
    51 = @{ E = "  +0.05%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
